# The document currently has two paragraphs:
#   1) "Лабораторная работа №3"
#   2) (empty, but holds the "_GoBack" bookmark)
#
# Target state:
#   1) same text, but now wrapped by the "_GoBack" bookmark
#   2) new heading text "Статистическое изучение взаимосвязи"
#      (centered, Times New Roman 14pt/ru-RU), bookmark removed from here.

$d = $word.ActiveDocument

# --- Move the hidden "_GoBack" bookmark onto paragraph 1's text -----------

# "_GoBack" is addressable by name even though it is hidden (starts with
# "_") and therefore excluded from Bookmarks.Count / enumeration.
$d.Bookmarks("_GoBack").Delete()

$p1 = $d.Paragraphs(1)
# Span just the run text (Start .. End-1), excluding the paragraph mark.
$bmRange = $d.Range($p1.Range.Start, $p1.Range.End - 1)
$d.Bookmarks.Add("_GoBack", $bmRange)

# --- Turn the (now bookmark-free) second paragraph into the new heading ---

$p2 = $d.Paragraphs(2)
$p2Range = $d.Range($p2.Range.Start, $p2.Range.End)

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:rFonts w:hint="default" w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="ru-RU"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="default" w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="ru-RU"/></w:rPr><w:t>Статистическое изучение взаимосвязи</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# InsertXML gives full control over the produced OOXML (formatting +
# w:hint="default"), unlike setting .Text/.Font piecemeal.
[void]$p2Range.InsertXML($xml)

# InsertXML inserted a brand-new paragraph ahead of the range's own trailing
# paragraph mark, leaving a stray empty paragraph behind it. Merge that
# leftover mark away so the document still has exactly two paragraphs.
$newP2 = $d.Paragraphs(2)
$newP3 = $d.Paragraphs(3)
$d.Range($newP2.Range.End - 1, $newP3.Range.End).Delete()
